# bug fixes on repost test 5
# Update column Q (Write Latency average) values for rows 3-23

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    3  = "8686.07"
    4  = "14938.35"
    5  = "74929.02"
    6  = "11302.20"
    7  = "55184.95"
    8  = "1297209.56"
    9  = "88738.80"
    10 = "8984.37"
    11 = "3842.78"
    12 = "15.82"
    13 = "9379.13"
    14 = "20.04"
    15 = "1357253.92"
    16 = "193189.34"
    17 = "12287.02"
    18 = "57778.46"
    19 = "323333.12"
    20 = "26.22"
    21 = "1844.93"
    22 = "79003.52"
    23 = "17702.64"
}

# Cells hold these figures as text (not numbers), so a leading apostrophe
# is used to force text entry and avoid numeric auto-conversion, matching
# the original cell's text storage.
foreach ($row in $values.Keys) {
    $cell = $ws.Range("Q$row")
    $cell.Value = "'" + $values[$row]
}
